# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Cebollín"
# right above the existing row 496. Every row from the old 496 through the old
# last row (537) shifts down by one (to 497..538), and the new row inherits the
# usual fixed fields for this market/category while carrying a new date and
# price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 496..537 down to 497..538, leaving a blank row 496 behind.
$ws.Rows(496).Insert()

# Populate the newly-inserted row 496 with the new observation.
$ws.Cells(496, 1).Value = 4
$ws.Cells(496, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells(496, 3).Value = "Los Lagos"
$ws.Cells(496, 4).Value = 45223
$ws.Cells(496, 5).Value = 10
$ws.Cells(496, 6).Value = 100112037
$ws.Cells(496, 7).Value = "Cebollín"
$ws.Cells(496, 8).Value = "Sin especificar"
$ws.Cells(496, 9).Value = "Primera"
$ws.Cells(496, 10).Value = 180
$ws.Cells(496, 11).Value = 6500
$ws.Cells(496, 12).Value = 6500
$ws.Cells(496, 13).Value = 6500
$ws.Cells(496, 14).Value = "$/paquete 36 unidades"
$ws.Cells(496, 15).Value = "Región Metropolitana"
$ws.Cells(496, 16).Value = 181
$ws.Cells(496, 17).Value = 36
$ws.Cells(496, 18).Value = "Hortaliza"
